$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update working_capital (column C) values
$ws.Range("C2").Value = 1000
$ws.Range("C3").Value = 1500
$ws.Range("C4").Value = 2000
$ws.Range("C5").Value = 3000
$ws.Range("C6").Value = 3500
$ws.Range("C7").Value = 4000
$ws.Range("C8").Value = 5000
$ws.Range("C9").Value = 5500
$ws.Range("C10").Value = 6000

# Update fixed_cost (column K) values
$ws.Range("K2").Value = 62
$ws.Range("K3").Value = 62
$ws.Range("K4").Value = 62
$ws.Range("K5").Value = 29
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 29
$ws.Range("K8").Value = 8
$ws.Range("K9").Value = 8
$ws.Range("K10").Value = 8

# Update the selection shown in sheetView
$ws.Range("K2:K4").Select()
